$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New recepcionista "Constanza" with zeros for prior weeks
$ws.Range("A6").Value = "Constanza"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# New column I header: week of 18_02_2024
$ws.Range("I1").Value = "18_02_2024"

# Existing recepcionistas get their value for the new week
$ws.Range("I2").Value = 935
$ws.Range("I3").Value = 902
$ws.Range("I4").Value = 1408
$ws.Range("I5").Value = 2904
$ws.Range("I6").Value = 20

$ws.Range("I7").Select()
